$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newK = @(3,5,5,5,4,3,7,2,1,3,5,2,0,0,7,4,7,7,8,4,7,4,4,1,5,3,6,7,5,4,4,2,0,0)

for ($i = 0; $i -lt $newK.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 7).Value = $newK[$i]
}
